$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# "Enterprises density (per 1000 people)" row: update Micro/SMEs/MSMEs values.
# These values are stored as text in the sheet (not numbers), so the leading
# apostrophe forces a text entry instead of letting Excel auto-convert the
# numeric-looking string to a Number. Re-applying the "Normal" style afterwards
# clears the transient quote-prefix/text-format styling so the cells keep
# their original (default/general) style, matching the source formatting.
$ws.Range("B13").Value = "'32.84"
$ws.Range("C13").Value = "'8.66"
$ws.Range("D13").Value = "'41.51"
$ws.Range("B13:D13").Style = "Normal"
